# PowerShell-style Word COM-interop script applying the SPT.docx edit.
# NOTE: text that contains a literal "$" (e.g. "${n#1}") MUST be written
# inside SINGLE quotes, otherwise PowerShell tries to expand it as a
# variable reference (${n#1} looks like ${varname} to the parser).

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Simple placeholder -> literal text replacements
# ------------------------------------------------------------------

# "${n#1}" -> "1"   (item number in the "Dasar" table)
$null = $d.Content.Find.Execute('${n#1}', $true, $false, $false, $false, $false, $true, 1, $false, '1', 2)

# "${i#1}" -> "1"   (item number in the "Kepada" table)
$null = $d.Content.Find.Execute('${i#1}', $true, $false, $false, $false, $false, $true, 1, $false, '1', 2)

# "${NAMA_PEGAWAI#1}" -> employee name
$null = $d.Content.Find.Execute('${NAMA_PEGAWAI#1}', $true, $false, $false, $false, $false, $true, 1, $false, 'ANDJAR SURJADIANTO,  S.Sos', 2)

# "melaksanakan" -> longer sentence (still starts with the same word)
$null = $d.Content.Find.Execute('melaksanakan', $true, $false, $false, $false, $false, $true, 1, $false, "melaksanakan Monitoring Pengelolaan Keuangan Desa di desa wilayah`nKecamatan Balongbendo Kabupaten Sidoarjo.", 2)

# "Jangka waktu" -> longer sentence
$null = $d.Content.Find.Execute('Jangka waktu', $true, $false, $false, $false, $false, $true, 1, $false, "Jangka waktu monitoring selama 8 (delapan) hari kerja pada periode tanggal 26`nJanuari s.d 4 Februari 2021.", 2)

# ------------------------------------------------------------------
# 2) "Dasar" table (table 1): replace uraian text and append 3 rows
#    describing the additional legal basis items (2, 3, 4).
# ------------------------------------------------------------------

$tDasar = $d.Tables.Item(1)

# "${uraian_dasar#1}" -> first legal-basis paragraph
$null = $d.Content.Find.Execute('${uraian_dasar#1}', $true, $false, $false, $false, $false, $true, 1, $false, "Undang - Undang Republik Indonesia Nomor 9 Tahun 2015 tentang`nPerubahan Kedua atas Undang-Undang Nomor 23 Tahun 2014`nPemerintah Daerah;", 2)

$dasarRow2 = $tDasar.Rows.Add()
$dasarRow3 = $tDasar.Rows.Add()
$dasarRow4 = $tDasar.Rows.Add()

$tDasar.Rows.Item(2).Cells.Item(2).Range.Text = '2.'
$tDasar.Rows.Item(2).Cells.Item(3).Range.Text = "Peraturan Pemerintah Republik Indonesia Nomor 12 Tahun 2017`ntentang Pembinaan dan Pengawasan Penyelenggaraan`nPemerintahan Daerah; "

$tDasar.Rows.Item(3).Cells.Item(2).Range.Text = '3.'
$tDasar.Rows.Item(3).Cells.Item(3).Range.Text = "Peraturan Menteri Dalam Negeri Nomor 23 Tahun 2020 tentang`nPerencanaan Pembinaan dan Pengawasan Pemerintahan Daerah`nTahun 2021; "

$tDasar.Rows.Item(4).Cells.Item(2).Range.Text = '4.'
$tDasar.Rows.Item(4).Cells.Item(3).Range.Text = "Program Kerja Pengawasan Tahunan (PKPT) Inspektorat Daerah`nKabupaten Sidoarjo Tahun 2021; "

# ------------------------------------------------------------------
# 3) "Kepada" table (table 2): replace the assignment/task text and
#    append 6 rows listing the rest of the assigned team members.
# ------------------------------------------------------------------

$tKepada = $d.Tables.Item(2)

# "${NAMA_TUGAS#1}" -> role of the first person
$null = $d.Content.Find.Execute('${NAMA_TUGAS#1}', $true, $false, $false, $false, $false, $true, 1, $false, 'Penanggungjawab', 2)

$kepadaRow2 = $tKepada.Rows.Add()
$kepadaRow3 = $tKepada.Rows.Add()
$kepadaRow4 = $tKepada.Rows.Add()
$kepadaRow5 = $tKepada.Rows.Add()
$kepadaRow6 = $tKepada.Rows.Add()
$kepadaRow7 = $tKepada.Rows.Add()

$tKepada.Rows.Item(2).Cells.Item(2).Range.Text = '2.'
$tKepada.Rows.Item(2).Cells.Item(3).Range.Text = 'Drs. PUJOSENO, SIP'
$tKepada.Rows.Item(2).Cells.Item(4).Range.Text = 'Pembantu Penanggungjawab'

$tKepada.Rows.Item(3).Cells.Item(2).Range.Text = '3.'
$tKepada.Rows.Item(3).Cells.Item(3).Range.Text = 'Drs. DANIEL S. TODING, M.AP'
$tKepada.Rows.Item(3).Cells.Item(4).Range.Text = 'Pengendali Mutu'

$tKepada.Rows.Item(4).Cells.Item(2).Range.Text = '4.'
$tKepada.Rows.Item(4).Cells.Item(3).Range.Text = 'HARI SUSANTO, S.Sos, MM'
$tKepada.Rows.Item(4).Cells.Item(4).Range.Text = 'Pengendali Teknis'

$tKepada.Rows.Item(5).Cells.Item(2).Range.Text = '5.'
$tKepada.Rows.Item(5).Cells.Item(3).Range.Text = 'MOCH. ARTFIANDO, SH'
$tKepada.Rows.Item(5).Cells.Item(4).Range.Text = 'Ketua Tim'

$tKepada.Rows.Item(6).Cells.Item(2).Range.Text = '6.'
$tKepada.Rows.Item(6).Cells.Item(3).Range.Text = 'NABILLAH CITRA CHAESARI, A.Md'
$tKepada.Rows.Item(6).Cells.Item(4).Range.Text = 'Anggota'

$tKepada.Rows.Item(7).Cells.Item(2).Range.Text = '7.'
$tKepada.Rows.Item(7).Cells.Item(3).Range.Text = 'YUANITA HASTOWO, A.Md'
$tKepada.Rows.Item(7).Cells.Item(4).Range.Text = 'Anggota'

Write-Output 'Edit complete'
